$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing four rows with the new text
$ws.Range("A1").Value = "Check Test: STARTED"
$ws.Range("A2").Value = "Check Test Case: startBrowser Test Method: SUCCESS"
$ws.Range("A3").Value = "Check Test Case: test Test Method: SUCCESS"
$ws.Range("A4").Value = "Check Test: ENDED"

# Remove the now-unused rows 5 and 6
$ws.Range("A5:A6").EntireRow.Delete()
